# Fixes #81 - Implemented the latest Excel import template.
#
# Rebuilds sheet1 from the old 14-column (A:N) layout into the new
# 19-column (A:S) layout: a "file path" column is inserted after
# "title", "agent:contributor" / "subject:spatial" columns are inserted
# after "agent:creator", "language" moves to just before "type of
# resource", and three new trailing columns are appended
# ("copyright jurisdiction", "copyright status").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Stash the handful of distinct cell formats already present in the
#    workbook into a holding row (far below the used range) so they
#    survive us overwriting the cells they currently live in. Re-using
#    these via Copy + PasteSpecial(Formats) keeps the existing font /
#    style table entries (and exact fractional sizes) instead of
#    minting new ones through the lossy Font.Size setter.
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("AA1").PasteSpecial($xlPasteFormats)     # bold 13.2 black  -> style 1
$ws.Range("C1").Copy()
$ws.Range("AA2").PasteSpecial($xlPasteFormats)     # bold 12 theme    -> style 2
$ws.Range("G1").Copy()
$ws.Range("AA3").PasteSpecial($xlPasteFormats)     # bold 12 black    -> style 3
$ws.Range("K2").Copy()
$ws.Range("AA4").PasteSpecial($xlPasteFormats)     # date numfmt      -> style 4
$ws.Range("B2").Copy()
$ws.Range("AA5").PasteSpecial($xlPasteFormats)     # right aligned    -> style 5
$ws.Range("D2").Copy()
$ws.Range("AA6").PasteSpecial($xlPasteFormats)     # plain 12 black   -> style 6
$ws.Range("A1").Copy()
$ws.Range("AA7").PasteSpecial($xlPasteFormats)     # plain 13.2 black -> style 7 (new)
$ws.Range("AA7").Font.Bold = $false

$excel.CutCopyMode = $false

$styleHeader = $ws.Range("AA1")
$styleBold12 = $ws.Range("AA2")
$styleBold12Black = $ws.Range("AA3")
$styleDate = $ws.Range("AA4")
$styleRight = $ws.Range("AA5")
$stylePlain12Black = $ws.Range("AA6")
$stylePlain132Black = $ws.Range("AA7")

# ---------------------------------------------------------------------
# 2. Clear the old A1:N4 grid completely (values + formats) so no stale
#    formatting lingers on cells that move or disappear.
# ---------------------------------------------------------------------
$ws.Range("A1:S4").Clear()

# ---------------------------------------------------------------------
# 3. Row 1 (headers)
# ---------------------------------------------------------------------
$row1 = @{
    "A1" = "object unique id";
    "B1" = "level";
    "C1" = "title";
    "D1" = "file path";
    "E1" = "file 1 name";
    "F1" = "file 1 use";
    "G1" = "file 2 name";
    "H1" = "file 2 name";
    "I1" = "subject:topic";
    "J1" = "agent:creator";
    "K1" = "agent:contributor";
    "L1" = "subject:spatial";
    "M1" = "date:created";
    "N1" = "note:note";
    "O1" = "identifier:doi";
    "P1" = "language";
    "Q1" = "type of resource";
    "R1" = "copyright jurisdiction";
    "S1" = "copyright status";
}
foreach ($addr in $row1.Keys) {
    $ws.Range($addr).Value = $row1[$addr]
}

$row1Style3 = @("H1", "K1", "R1")
$row1Style1 = @("A1", "B1", "Q1", "S1")
foreach ($addr in $row1.Keys) {
    if ($row1Style3 -contains $addr) {
        $styleBold12Black.Copy()
    } elseif ($row1Style1 -contains $addr) {
        $styleHeader.Copy()
    } else {
        $styleBold12.Copy()
    }
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Row 2 (object-level sample data)
# ---------------------------------------------------------------------
$row2 = @{
    "A2" = "object#1";
    "B2" = "Object";
    "C2" = "Test Object One";
    "I2" = "SUBJECT:TOPIC";
    "J2" = "AGENT:CREATOR";
    "K2" = "AGENT:CONTRIBUTOR";
    "L2" = "SUBJECT:SPATIAL";
    "M2" = "Decenber 10, 2016 @{begin=2016-12-10 ; end=2016-12-10 }";
    "N2" = "NOTE:NOTE";
    "O2" = "IDENTIFIER:DOI";
    "P2" = "eng  - English|zxx  - No linguistic content; Not applicable ";
    "Q2" = "mixed material | still image";
    "R2" = "US - United States of America";
    "S2" = "copyrighted";
}
foreach ($addr in $row2.Keys) {
    $ws.Range($addr).Value = $row2[$addr]
}
# E2 / F2 stay empty but carry the plain-black style, same as before.
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""

$row2StylePlain12 = @("E2", "F2", "K2", "P2", "Q2")
$row2StylePlain132 = @("R2", "S2")

foreach ($addr in $row2StylePlain12) {
    $stylePlain12Black.Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}
foreach ($addr in $row2StylePlain132) {
    $stylePlain132Black.Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}
$styleRight.Copy()
$ws.Range("B2").PasteSpecial($xlPasteFormats)
$styleDate.Copy()
$ws.Range("M2").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 5. Rows 3 & 4 (component / sub-component sample data) - plain cells.
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "object#1"
$ws.Range("B3").Value = "Component"
$ws.Range("C3").Value = "Test Component One"
$ws.Range("E3").Value = "file_1.jpg"
$ws.Range("F3").Value = "image-source"

$ws.Range("A4").Value = "object#1"
$ws.Range("B4").Value = "Sub-component"
$ws.Range("C4").Value = "Test Sub-component One"
$ws.Range("E4").Value = "file_2.jpg"
$ws.Range("F4").Value = "image-source"

# ---------------------------------------------------------------------
# 6. Column widths.
# ---------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 14.6    # -> stored width 15.5
$ws.Columns("D").ColumnWidth = 14.6    # -> stored width 15.5
$ws.Columns("N").ColumnWidth = 38.75   # -> stored width ~39.664 (note:note)
$ws.Columns("O").ColumnWidth = 29.084  # -> stored width 30 (identifier:doi)
$ws.Columns("Q").ColumnWidth = 25.92   # -> stored width ~26.832 (type of resource)

# ---------------------------------------------------------------------
# 7. Selection / view state.
# ---------------------------------------------------------------------
$ws.Range("K1").Select()

# ---------------------------------------------------------------------
# 8. Clean up the scratch holding area used for format templates.
# ---------------------------------------------------------------------
$ws.Range("AA1:AA7").Clear()

$ws.Range("K1").Select()
